$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the existing header style (bold font, border, centered) to the new
# label column (A2:A6) and to the new header cell (C1), by copying format
# + value from the current "média/animal" header cell, then overwrite the
# values afterwards.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("B1").Copy($ws.Range("A2"))
$ws.Range("B1").Copy($ws.Range("A3"))
$ws.Range("B1").Copy($ws.Range("A4"))
$ws.Range("B1").Copy($ws.Range("A5"))
$ws.Range("B1").Copy($ws.Range("A6"))

# --- Header row: shift from A1:B1 to B1:C1, clear A1 ---
$ws.Range("C1").Value = "média/animal"
$ws.Range("B1").Value = "total"
$ws.Range("A1").Clear()

# --- Data rows: shift existing numeric columns A:B -> B:C ---
$ws.Range("B2").Value = 16068
$ws.Range("C2").Value = 16

$ws.Range("B3").Value = 15980
$ws.Range("C3").Value = 16

$ws.Range("B4").Value = 1864
$ws.Range("C4").Value = 2

$ws.Range("B5").Value = 6019
$ws.Range("C5").Value = 6

$ws.Range("B6").Value = 6414
$ws.Range("C6").Value = 6

# --- New label column A values ---
$ws.Range("A2").Value = "tratamentos sarna"
$ws.Range("A3").Value = "tratamentos pulga"
$ws.Range("A4").Value = "vacinas"
$ws.Range("A5").Value = "vermífugos"
$ws.Range("A6").Value = "dias de internação"
